# "new pp set and refinement zones"
# Clear the old single-point reference ("h_ref" label + its value/diff column)
# and blank out the refinement-zone rows (17-26), leaving only the empty,
# formatted D column cells that mark the new pp set placeholder rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C16 held the "h_ref" label, C17:C26 held the manual reference values,
# D17:D26 held the "=Cxx-Cyy" delta formulas. Clear all of column C (16:26)
# and the formulas in D17:D26 so only the blank, styled D cells remain.
$ws.Range("C16:C26").ClearContents()
$ws.Range("D17:D26").ClearContents()

# Update the view: scroll down to row 5 and select the newly cleared
# refinement-zone range C16:D26 with C16 as the active cell.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("C16:D26").Select() | Out-Null
